$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "MCT-2A-Eletrônica analóg. e de potência"

$ws.Range("E3").Value = "-"
$ws.Range("F3").Value = "-"

$ws.Range("B8").Value = "MCT-2A-Eletrônica analóg. e de potência"

$ws.Range("C11").Value = "[-, -, -, 'MEC-1A-Comandos Eletricos']"
$ws.Range("E11").Value = "-"
$ws.Range("F11").Value = "[-, -, -, 'MEC-2A-Elet. Dig. Bas.']"

$ws.Range("C12").Value = "[-, -, -, 'MEC-1A-Comandos Eletricos']"
$ws.Range("E12").Value = "-"
$ws.Range("F12").Value = "[-, -, -, 'MEC-2A-Elet. Dig. Bas.']"

$ws.Range("C14").Value = "[-, -, -, 'MEC-1A-Comandos Eletricos']"
$ws.Range("E14").Value = "-"
$ws.Range("F14").Value = "[-, -, -, 'MEC-2A-Elet. Dig. Bas.']"

$ws.Range("C15").Value = "[-, -, -, 'MEC-1A-Comandos Eletricos']"
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "[-, -, -, 'MEC-2A-Elet. Dig. Bas.']"

$ws.Range("D18").Value = "[-, -, 'MEC-1NB-Elet. Dig. Bas.', -]"

$ws.Range("C19").Value = "ELM-2NA-Circuitos Elétricos 2"
$ws.Range("D19").Value = "[-, -, 'MEC-1NB-Elet. Dig. Bas.', -]"
$ws.Range("F19").Value = "[-, -, 'MEC-1NB-Elet. Dig. Bas.', -]"

$ws.Range("D20").Value = "[-, -, -, 'MEC-1NB-Elet. Dig. Bas.']"
$ws.Range("F20").Value = "-"

$ws.Range("C21").Value = "ELM-2NA-Circuitos Elétricos 2"
$ws.Range("F21").Value = "-"
